$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = '{''shimenet'', ''$'', ''andamhie'', ''anda'', ''chika'', ''eklabool''}'
$ws.Range("D4").Value = '{''shimenet'', ''naur'', ''$'', ''andamhie'', ''anda'', ''chika'', ''eklabool''}'
$ws.Range("D5").Value = '{''chika'', ''andamhie'', ''anda'', ''eklabool''}'
$ws.Range("D11").Value = '{'')'', '',''}'
$ws.Range("D13").Value = '{'';'', '')'', '','', ''=''}'
$ws.Range("D14").Value = '{''/='', ''=='', ''**'', ''}'', '';'', ''/'', ''%'', ''*='', ''>'', ''step'', ''-='', ''<'', '')'', ''<='', '':'', ''to'', ''**='', ''//='', '']'', ''%='', ''*'', ''+='', ''+'', ''>='', ''//'', ''&&'', ''!='', '','', ''-'', ''='', ''||''}'
$ws.Range("D15").Value = '{''/='', ''=='', ''**'', ''}'', '';'', ''/'', ''%'', ''*='', ''>'', ''step'', ''-='', ''<'', '')'', ''<='', '':'', ''to'', ''**='', ''//='', '']'', ''%='', ''*'', ''+='', ''+'', ''>='', ''//'', ''&&'', ''!='', '','', ''-'', ''='', ''||''}'
$ws.Range("D16").Value = '{''/='', ''=='', '')'', ''<='', '':'', ''to'', ''**='', ''**'', ''//='', '']'', ''}'', ''%='', '';'', ''*'', ''+'', ''>='', ''/'', ''%'', ''*='', ''='', ''>'', ''//'', ''&&'', ''!='', '','', ''||'', ''step'', ''-'', ''-='', ''<'', ''+=''}'
$ws.Range("D17").Value = '{'';'', '','', '']'', ''}''}'
$ws.Range("D18").Value = '{''=='', ''**'', ''}'', '';'', ''/'', ''%'', ''>'', ''step'', ''<'', '')'', ''<='', '':'', ''to'', '']'', ''*'', ''+'', ''>='', ''//'', ''&&'', ''!='', '','', ''-'', ''||''}'
$ws.Range("D19").Value = '{''=='', '')'', ''<='', '':'', ''to'', ''**'', '']'', ''}'', '';'', ''*'', ''+'', ''>='', ''/'', ''%'', ''>'', ''//'', ''&&'', ''!='', '','', ''step'', ''-'', ''<'', ''||''}'
$ws.Range("D20").Value = '{''=='', ''**'', ''}'', '';'', ''/'', ''%'', ''>'', ''step'', ''<'', '')'', ''<='', '':'', ''to'', '']'', ''*'', ''+'', ''>='', ''//'', ''&&'', ''!='', '','', ''-'', ''||''}'
$ws.Range("D21").Value = '{''=='', ''**'', ''}'', '';'', ''/'', ''%'', ''>'', ''step'', ''<'', '')'', ''<='', '':'', ''to'', '']'', ''*'', ''+'', ''>='', ''//'', ''&&'', ''!='', '','', ''-'', ''||''}'
$ws.Range("D32").Value = '{''serve'', ''}'', ''eklabool'', ''keri'', ''versa'', ''gogogo'', ''id'', ''betsung'', ''pak'', ''--'', ''chika'', ''push'', ''++'', ''ditech'', ''naur'', ''adelete'', ''amaccana'', ''andamhie'', ''anda'', ''forda'', ''adele''}'
$ws.Range("D33").Value = '{''id'', ''betsung'', ''serve'', ''pak'', ''gogogo'', ''}'', ''chika'', ''eklabool'', ''--'', ''push'', ''++'', ''ditech'', ''naur'', ''adelete'', ''amaccana'', ''andamhie'', ''anda'', ''keri'', ''versa'', ''forda'', ''adele''}'
$ws.Range("D36").Value = '{'')'', '',''}'
$ws.Range("D37").Value = '{''serve'', ''}'', ''eklabool'', ''keri'', ''versa'', ''gogogo'', ''id'', ''betsung'', ''pak'', ''--'', ''chika'', ''push'', ''++'', ''ditech'', ''naur'', ''adelete'', ''amaccana'', ''andamhie'', ''anda'', ''forda'', ''adele''}'
$ws.Range("D39").Value = '{'';'', '')'', '':'', ''to'', '','', '']'', ''step'', ''}''}'
$ws.Range("D40").Value = '{'';'', '')'', '':'', ''to'', '','', '']'', ''step'', ''}''}'
$ws.Range("D41").Value = '{'';'', '')'', '':'', ''to'', '','', '']'', ''step'', ''}''}'
$ws.Range("D42").Value = '{'';'', '')'', '':'', ''to'', '','', '']'', ''step'', ''}''}'
$ws.Range("D43").Value = '{''=='', ''**'', ''}'', '';'', ''/'', ''%'', ''>'', ''step'', ''<'', '')'', ''<='', '':'', ''to'', '']'', ''*'', ''+'', ''>='', ''//'', ''&&'', ''!='', '','', ''-'', ''||''}'
$ws.Range("D44").Value = '{''andamhie_literal'', ''('', ''id'', ''len'', ''korik'', ''eme'', ''chika_literal'', ''--'', ''++'', ''anda_literal''}'
$ws.Range("D45").Value = '{''=='', ''**'', ''}'', '';'', ''/'', ''%'', ''>'', ''step'', ''<'', '')'', ''<='', '':'', ''to'', '']'', ''*'', ''+'', ''>='', ''//'', ''&&'', ''!='', '','', ''-'', ''||''}'
$ws.Range("D47").Value = '{''=='', ''**'', ''}'', '';'', ''/'', ''%'', ''>'', ''step'', ''<'', '')'', ''<='', '':'', ''to'', '']'', ''*'', ''+'', ''>='', ''//'', ''&&'', ''!='', '','', ''-'', ''||''}'
$ws.Range("D48").Value = '{''id'', ''=='', '')'', ''<='', '':'', ''to'', ''**'', '']'', ''}'', '';'', ''*'', ''+'', ''>='', ''/'', ''%'', ''>'', ''//'', ''&&'', ''!='', '','', ''step'', ''-'', ''<'', ''||''}'
$ws.Range("D49").Value = '{''=='', ''**'', ''}'', '';'', ''/'', ''%'', ''>'', ''step'', ''<'', '')'', ''<='', '':'', ''to'', '']'', ''*'', ''+'', ''>='', ''//'', ''&&'', ''!='', '','', ''-'', ''||''}'
$ws.Range("D50").Value = '{''andamhie_literal'', ''('', ''id'', ''len'', ''korik'', ''eme'', ''--'', ''chika_literal'', ''-'', ''!'', ''++'', ''anda_literal''}'
$ws.Range("D51").Value = '{''forda'', ''}'', ''eklabool'', ''keri'', ''versa'', ''gogogo'', ''id'', ''betsung'', ''pak'', ''--'', ''chika'', ''push'', ''++'', ''ditech'', ''naur'', ''adelete'', ''amaccana'', ''andamhie'', ''anda'', ''serve'', ''adele''}'
$ws.Range("D52").Value = '{''forda'', ''}'', ''eklabool'', ''keri'', ''versa'', ''gogogo'', ''id'', ''betsung'', ''pak'', ''--'', ''chika'', ''push'', ''++'', ''ditech'', ''naur'', ''adelete'', ''amaccana'', ''andamhie'', ''anda'', ''serve'', ''adele''}'
$ws.Range("D53").Value = '{''betsung'', ''}'', ''ditech''}'
$ws.Range("D54").Value = '{''serve'', ''}'', ''eklabool'', ''keri'', ''versa'', ''gogogo'', ''id'', ''betsung'', ''pak'', ''--'', ''chika'', ''push'', ''++'', ''ditech'', ''naur'', ''adelete'', ''amaccana'', ''andamhie'', ''anda'', ''forda'', ''adele''}'
$ws.Range("D55").Value = '{''andamhie_literal'', ''('', ''id'', ''len'', ''korik'', ''eme'', ''{'', ''--'', ''chika_literal'', ''-'', ''!'', ''++'', ''anda_literal''}'
$ws.Range("D57").Value = '{''serve'', ''}'', ''eklabool'', ''keri'', ''versa'', ''gogogo'', ''id'', ''betsung'', ''pak'', ''--'', ''chika'', ''push'', ''++'', ''ditech'', ''naur'', ''adelete'', ''amaccana'', ''andamhie'', ''anda'', ''forda'', ''adele''}'
$ws.Range("D60").Value = '{''serve'', ''}'', ''eklabool'', ''keri'', ''versa'', ''gogogo'', ''id'', ''betsung'', ''pak'', ''--'', ''chika'', ''push'', ''++'', ''ditech'', ''naur'', ''adelete'', ''amaccana'', ''andamhie'', ''anda'', ''forda'', ''adele''}'
$ws.Range("D61").Value = '{''serve'', ''}'', ''eklabool'', ''keri'', ''versa'', ''gogogo'', ''id'', ''betsung'', ''pak'', ''--'', ''chika'', ''push'', ''++'', ''ditech'', ''naur'', ''adelete'', ''amaccana'', ''andamhie'', ''anda'', ''forda'', ''adele''}'
$ws.Range("D62").Value = '{''serve'', ''}'', ''eklabool'', ''keri'', ''versa'', ''gogogo'', ''id'', ''betsung'', ''pak'', ''--'', ''chika'', ''push'', ''++'', ''ditech'', ''naur'', ''adelete'', ''amaccana'', ''andamhie'', ''anda'', ''forda'', ''adele''}'
$ws.Range("D64").Value = '{''serve'', ''}'', ''eklabool'', ''keri'', ''versa'', ''gogogo'', ''id'', ''betsung'', ''pak'', ''--'', ''chika'', ''push'', ''++'', ''ditech'', ''naur'', ''adelete'', ''amaccana'', ''andamhie'', ''anda'', ''forda'', ''adele''}'
$ws.Range("D67").Value = '{''serve'', ''}'', ''eklabool'', ''keri'', ''versa'', ''gogogo'', ''id'', ''betsung'', ''pak'', ''--'', ''chika'', ''push'', ''++'', ''ditech'', ''naur'', ''adelete'', ''amaccana'', ''andamhie'', ''anda'', ''forda'', ''adele''}'
$ws.Range("D69").Value = '{''serve'', ''ganern'', ''}'', ''eklabool'', ''keri'', ''versa'', ''gogogo'', ''id'', ''betsung'', ''pak'', ''--'', ''chika'', ''push'', ''++'', ''ditech'', ''naur'', ''adelete'', ''amaccana'', ''andamhie'', ''anda'', ''forda'', ''adele''}'
$ws.Range("D70").Value = '{''serve'', ''}'', ''eklabool'', ''keri'', ''versa'', ''gogogo'', ''id'', ''betsung'', ''pak'', ''--'', ''chika'', ''push'', ''++'', ''ditech'', ''naur'', ''adelete'', ''amaccana'', ''andamhie'', ''anda'', ''forda'', ''adele''}'
$ws.Range("D71").Value = '{''serve'', ''}'', ''eklabool'', ''keri'', ''versa'', ''gogogo'', ''id'', ''betsung'', ''pak'', ''--'', ''chika'', ''push'', ''++'', ''ditech'', ''naur'', ''adelete'', ''amaccana'', ''andamhie'', ''anda'', ''forda'', ''adele''}'
$ws.Range("D72").Value = '{''forda'', ''id'', ''betsung'', ''pak'', ''gogogo'', ''}'', ''eklabool'', ''--'', ''chika'', ''push'', ''++'', ''ditech'', ''naur'', ''adelete'', ''amaccana'', ''andamhie'', ''anda'', ''keri'', ''versa'', ''serve'', ''adele''}'
$ws.Range("D78").Value = '{'')'', ''step'', ''to''}'
$ws.Range("D80").Value = '{''serve'', ''}'', ''eklabool'', ''keri'', ''versa'', ''gogogo'', ''id'', ''betsung'', ''pak'', ''--'', ''chika'', ''push'', ''++'', ''ditech'', ''naur'', ''adelete'', ''amaccana'', ''andamhie'', ''anda'', ''forda'', ''adele''}'
$ws.Range("D81").Value = '{''forda'', ''id'', ''betsung'', ''pak'', ''gogogo'', ''}'', ''eklabool'', ''--'', ''chika'', ''push'', ''++'', ''ditech'', ''naur'', ''adelete'', ''amaccana'', ''andamhie'', ''anda'', ''keri'', ''versa'', ''serve'', ''adele''}'
$ws.Range("D83").Value = '{''forda'', ''id'', ''betsung'', ''pak'', ''gogogo'', ''--'', ''chika'', ''push'', ''eklabool'', ''++'', ''}'', ''ditech'', ''naur'', ''adelete'', ''amaccana'', ''andamhie'', ''anda'', ''keri'', ''versa'', ''serve'', ''adele''}'
$ws.Range("D86").Value = '{''amaccana'', ''betsung'', ''}'', ''ditech''}'
$ws.Range("D87").Value = '{''serve'', ''}'', ''eklabool'', ''keri'', ''versa'', ''gogogo'', ''id'', ''betsung'', ''pak'', ''--'', ''chika'', ''push'', ''++'', ''ditech'', ''naur'', ''adelete'', ''amaccana'', ''andamhie'', ''anda'', ''forda'', ''adele''}'
$ws.Range("D89").Value = '{''betsung'', ''}'', ''ditech''}'
$ws.Range("D90").Value = '{''betsung'', ''}'', ''ditech''}'
$ws.Range("D92").Value = '{''serve'', ''}'', ''eklabool'', ''keri'', ''versa'', ''gogogo'', ''id'', ''betsung'', ''pak'', ''--'', ''chika'', ''push'', ''++'', ''ditech'', ''naur'', ''adelete'', ''amaccana'', ''andamhie'', ''anda'', ''forda'', ''adele''}'
$ws.Range("D93").Value = '{''serve'', ''}'', ''eklabool'', ''keri'', ''versa'', ''gogogo'', ''id'', ''betsung'', ''pak'', ''--'', ''chika'', ''push'', ''++'', ''ditech'', ''naur'', ''adelete'', ''amaccana'', ''andamhie'', ''anda'', ''forda'', ''adele''}'
$ws.Range("D94").Value = '{''serve'', ''}'', ''eklabool'', ''keri'', ''versa'', ''gogogo'', ''id'', ''betsung'', ''pak'', ''--'', ''chika'', ''push'', ''++'', ''ditech'', ''naur'', ''adelete'', ''amaccana'', ''andamhie'', ''anda'', ''forda'', ''adele''}'
